$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the data of row 2 and row 3 for the specific subset of
# columns whose values actually differ between the two records:
#   A, I, P, Q, R, S, Y, AA, AC, AW, AX
# All other columns hold identical values in both rows, so they are left
# untouched.

# Plain columns: safe to copy via Value2 without Excel trying to
# "smart convert" the text into something else.
$plainCols = @("A", "P", "Q", "R", "S", "AC", "AW", "AX")

# Text columns: Excel would otherwise auto-convert these into a number
# (I holds numeric-looking text like "1") or a date serial (Y/AA hold
# date-like text such as "2022-08-08"). Force Text format first so the
# values are kept as literal text, then restore the Normal style so no
# stray cell formatting is left behind.
$textCols = @("I", "Y", "AA")

foreach ($col in $plainCols) {
    $addr2 = $col + "2"
    $addr3 = $col + "3"

    $val2 = $ws.Range($addr2).Value2
    $val3 = $ws.Range($addr3).Value2

    if ($null -eq $val2) { $val2 = "" }
    if ($null -eq $val3) { $val3 = "" }

    $ws.Range($addr2).Value2 = $val3
    $ws.Range($addr3).Value2 = $val2
}

foreach ($col in $textCols) {
    $addr2 = $col + "2"
    $addr3 = $col + "3"

    $val2 = $ws.Range($addr2).Value2
    $val3 = $ws.Range($addr3).Value2

    if ($null -eq $val2) { $val2 = "" }
    if ($null -eq $val3) { $val3 = "" }

    $ws.Range($addr2).NumberFormat = "@"
    $ws.Range($addr3).NumberFormat = "@"

    $ws.Range($addr2).Value2 = $val3
    $ws.Range($addr3).Value2 = $val2

    $ws.Range($addr2).Style = "Normal"
    $ws.Range($addr3).Style = "Normal"
}
